$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$participantQuery = @'
WITH {
    phs_accession: "phs002517",
    subject_ids: [],
    experimental_strategies: ["WGS"],
    genders: ["Female"],
    sample_tumor_statuses: [],
    file_types: [],
    library_strategies: [],
    library_sources: [],
    library_selections: [],
    library_layouts: [],
    platforms: [],
    instrument_models: [],
    reference_genome_assemblies: [],
    primary_diagnoses: [],
    num_study_samples_min: 0,
    num_study_samples_max: 0,
    num_study_participants_max: 0,
    num_study_participants_min: 0
} AS inputs, "Not specified in data" AS na
MATCH (s:study {phs_accession: inputs.phs_accession})
OPTIONAL MATCH (s)<--(p:participant)
WITH inputs, na, s, count(distinct p) AS num_p
WHERE
    (inputs.num_study_participants_min = 0 OR num_p >= inputs.num_study_participants_min) AND
    (inputs.num_study_participants_max = 0 OR num_p <= inputs.num_study_participants_max)
OPTIONAL MATCH (s)<--(:participant)<--(samp:sample)
WITH inputs, na, s, count(distinct samp) AS num_samp
WHERE
    (inputs.num_study_samples_min = 0 OR num_samp >= inputs.num_study_samples_min) AND
    (inputs.num_study_participants_max = 0 OR num_samp <= inputs.num_study_participants_max)
MATCH (s)<--(p:participant)
WITH inputs, na, p, {
	participant_id: coalesce(p.participant_id, ""),
	study_name: coalesce(s.study_name, ""),
	accession: coalesce(s.phs_accession, ""),
	gender: coalesce(p.gender, "")
} AS data
WHERE 
    (size(inputs.subject_ids) = 0 OR p.participant_id IN inputs.subject_ids) AND
    (size(inputs.genders) = 0 OR p.gender IN inputs.genders)
OPTIONAL MATCH (p)<--(samp:sample)
WITH inputs, na, p, apoc.map.merge(data, {
		samples: coalesce(apoc.text.join(apoc.coll.sort(collect(distinct samp.sample_id)), ""), "")
	}) AS data,
	COLLECT(DISTINCT samp.sample_tumor_status) AS sample_tumor_statuses
WHERE
    (size(inputs.sample_tumor_statuses) = 0 OR size(apoc.coll.intersection(inputs.sample_tumor_statuses, sample_tumor_statuses)) > 0)
OPTIONAL MATCH (p)<--(:sample)<--(f:file)
WITH inputs, na, p, data,
	COLLECT(DISTINCT f.file_type) AS file_types, 
	apoc.coll.toSet(apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\\s{0,1}")) AS experimental_strategies
WHERE
    (size(inputs.experimental_strategies) = 0 OR size(apoc.coll.intersection(inputs.experimental_strategies, experimental_strategies)) > 0) AND
    (size(inputs.file_types) = 0 OR size(apoc.coll.intersection(inputs.file_types, file_types)) > 0)
OPTIONAL MATCH (p)<--(:sample)<--(:file)<--(g:genomic_info)
WITH inputs, na, p, data,
	COLLECT(DISTINCT g.library_strategy) AS library_strategy,
	COLLECT(DISTINCT g.library_source) AS library_source,
	COLLECT(DISTINCT g.library_selection) AS library_selection,
	COLLECT(DISTINCT g.library_layout) AS library_layout,
	COLLECT(DISTINCT g.platform) AS platform,
	COLLECT(DISTINCT g.instrument_model) AS instrument_models,
	COLLECT(DISTINCT g.reference_genome_assembly) AS reference_genome_assembly
WHERE
    (size(inputs.library_strategies) = 0 OR size(apoc.coll.intersection(inputs.library_strategies, library_strategy)) > 0) AND
    (size(inputs.library_sources) = 0 OR size(apoc.coll.intersection(inputs.library_sources, library_source)) > 0) AND 
    (size(inputs.library_selections) = 0 OR size(apoc.coll.intersection(inputs.library_selections, library_selection)) > 0) AND
    (size(inputs.library_layouts) = 0 OR size(apoc.coll.intersection(inputs.library_layouts, library_layout)) > 0) AND
    (size(inputs.platforms) = 0 OR size(apoc.coll.intersection(inputs.platforms, platform)) > 0) AND
    (size(inputs.instrument_models) = 0 OR size(apoc.coll.intersection(inputs.instrument_models, instrument_models)) > 0) AND
    (size(inputs.reference_genome_assemblies) = 0 OR size(apoc.coll.intersection(inputs.reference_genome_assemblies, reference_genome_assembly)) > 0)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH inputs, na, p, data,
	COLLECT(DISTINCT diag.primary_diagnosis) AS primary_diagnosis
WHERE
    (size(inputs.primary_diagnoses) = 0 OR size(apoc.coll.intersection(inputs.primary_diagnoses, primary_diagnosis)) > 0)
RETURN 
	data.participant_id AS `Participant ID`,
	data.study_name AS `Study Name`,
	data.accession AS `Accession`,
	data.gender AS `Gender`,
	data.samples AS `Samples`
ORDER BY `Participant ID`
LIMIT 100
'@

$sampleQuery = @'
WITH {
    phs_accession: "phs002517",
    subject_ids: [],
    experimental_strategies: ["WGS"],
    gender: ["Female"],
    sample_tumor_status: [],
    file_types: [],
    library_strategies: [],
    library_sources: [],
    library_selections: [],
    library_layouts: [],
    platforms: [],
    instrument_models: [],
    reference_genome_assemblies: [],
    primary_diagnoses: [],
    num_study_samples_min: 0,
    num_study_samples_max: 0,
    num_study_participants_max: 0,
    num_study_participants_min: 0
} AS inputs, "Not specified in data" AS na
MATCH (s:study {phs_accession: inputs.phs_accession})
OPTIONAL MATCH (s)<--(p:participant)
WITH inputs, na, s, count(distinct p) AS num_p
WHERE
    (inputs.num_study_participants_min = 0 OR num_p >= inputs.num_study_participants_min) AND
    (inputs.num_study_participants_max = 0 OR num_p <= inputs.num_study_participants_max)
OPTIONAL MATCH (s)<--(:participant)<--(samp:sample)
WITH inputs, na, s, count(distinct samp) AS num_samp
WHERE
    (inputs.num_study_samples_min = 0 OR num_samp >= inputs.num_study_samples_min) AND
    (inputs.num_study_participants_max = 0 OR num_samp <= inputs.num_study_participants_max)
MATCH (s)<--(:participant)<--(samp:sample)
WITH inputs, na, samp, {
	sample_id: samp.sample_id,
	study_name: s.study_name,
	accession: s.phs_accession,
	tumor: samp.sample_tumor_status,
	analyte_type: samp.sample_type
} AS data
WHERE
    (size(inputs.sample_tumor_status) = 0 OR samp.sample_tumor_status IN inputs.sample_tumor_status)
OPTIONAL MATCH (samp)-->(p:participant)
WITH inputs, na, samp, 
	apoc.map.merge(data, {
		participant_id: p.participant_id
	}) AS data,
	COLLECT(DISTINCT p.participant_id) AS participant_id,
	COLLECT(DISTINCT p.gender) AS gender
WHERE 
    (size(inputs.subject_ids) = 0 OR size(apoc.coll.intersection(inputs.participant_id, participant_id)) > 0) AND
    (size(inputs.gender) = 0 OR size(apoc.coll.intersection(inputs.gender, gender)) > 0)
OPTIONAL MATCH (samp)<--(f:file)
WITH inputs, na, samp, data,
	COLLECT(DISTINCT f.file_type) AS file_types, 
	apoc.coll.toSet(apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\\s{0,1}")) AS experimental_strategies
WHERE
	(size(inputs.experimental_strategies) = 0 OR size(apoc.coll.intersection(inputs.experimental_strategies, experimental_strategies)) > 0) AND
    (size(inputs.file_types) = 0 OR size(apoc.coll.intersection(inputs.file_types, file_types)) > 0)
OPTIONAL MATCH (samp)<--(:file)<--(g:genomic_info)
WITH inputs, na, samp, data,
	COLLECT(DISTINCT g.library_strategy) AS library_strategy,
	COLLECT(DISTINCT g.library_source) AS library_source,
	COLLECT(DISTINCT g.library_selection) AS library_selection,
	COLLECT(DISTINCT g.library_layout) AS library_layout,
	COLLECT(DISTINCT g.platform) AS platform,
	COLLECT(DISTINCT g.instrument_model) AS instrument_models,
	COLLECT(DISTINCT g.reference_genome_assembly) AS reference_genome_assembly
WHERE
    (size(inputs.library_strategies) = 0 OR size(apoc.coll.intersection(inputs.library_strategies, library_strategy)) > 0) AND
    (size(inputs.library_sources) = 0 OR size(apoc.coll.intersection(inputs.library_sources, library_source)) > 0) AND 
    (size(inputs.library_selections) = 0 OR size(apoc.coll.intersection(inputs.library_selections, library_selection)) > 0) AND
    (size(inputs.library_layouts) = 0 OR size(apoc.coll.intersection(inputs.library_layouts, library_layout)) > 0) AND
    (size(inputs.platforms) = 0 OR size(apoc.coll.intersection(inputs.platforms, platform)) > 0) AND
    (size(inputs.instrument_models) = 0 OR size(apoc.coll.intersection(inputs.instrument_models, instrument_models)) > 0) AND
    (size(inputs.reference_genome_assemblies) = 0 OR size(apoc.coll.intersection(inputs.reference_genome_assemblies, reference_genome_assembly)) > 0)
OPTIONAL MATCH (samp)-->(:participant)<--(diag:diagnosis)
WITH inputs, na, samp, data,
	COLLECT(DISTINCT diag.primary_diagnosis) AS primary_diagnosis
WHERE
    (size(inputs.primary_diagnoses) = 0 OR size(apoc.coll.intersection(inputs.primary_diagnoses, primary_diagnosis)) > 0)
RETURN 
	data.sample_id AS `Sample ID`,
	data.participant_id AS `Participant ID`,
	data.study_name AS `Study Name`,
	data.accession AS `Accession`,
	data.tumor AS `Tumor`,
	data.analyte_type AS `Analyte Type`
ORDER BY `Sample ID`
LIMIT 100
'@

$fileQuery = @'
WITH {
    phs_accession: "phs002517",
    subject_ids: [],
    experimental_strategies: ["WGS"],
    genders: ["Female"],
    sample_tumor_statuses: [],
    file_types: [],
    library_strategies: [],
    library_sources: [],
    library_selections: [],
    library_layouts: [],
    platforms: [],
    instrument_models: [],
    reference_genome_assemblies: [],
    primary_diagnoses: [],
    num_study_samples_min: 0,
    num_study_samples_max: 0,
    num_study_participants_max: 0,
    num_study_participants_min: 0
} AS inputs, "Not specified in data" AS na
MATCH (s:study {phs_accession: inputs.phs_accession})
OPTIONAL MATCH (s)<--(p:participant)
WITH inputs, na, s, count(distinct p) AS num_p
WHERE
    (inputs.num_study_participants_min = 0 OR num_p >= inputs.num_study_participants_min) AND
    (inputs.num_study_participants_max = 0 OR num_p <= inputs.num_study_participants_max)
OPTIONAL MATCH (s)<--(:participant)<--(samp:sample)
WITH inputs, na, s, count(distinct samp) AS num_samp
WHERE
    (inputs.num_study_samples_min = 0 OR num_samp >= inputs.num_study_samples_min) AND
    (inputs.num_study_participants_max = 0 OR num_samp <= inputs.num_study_participants_max)
MATCH (s)<--(f:file)
WITH inputs, na, f, {
        file_name: coalesce(f.file_name, ""),
        file_type: coalesce(f.file_type, ""),
        study_name: coalesce(s.study_name, ""),
        accession: coalesce(s.phs_accession, "")
    } AS data,
    apoc.coll.toSet(apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\\s{0,1}")) AS experimental_strategies
WHERE 
    (size(inputs.experimental_strategies) = 0 OR size(apoc.coll.intersection(inputs.experimental_strategies, experimental_strategies)) > 0) AND
    (size(inputs.file_types) = 0 OR f.file_type IN inputs.file_types)
OPTIONAL MATCH (f)-->(samp)
WITH inputs, na, f, apoc.map.merge(data, {
        sample_id: coalesce(apoc.text.join(apoc.coll.sort(collect(distinct samp.sample_id)), ", "), "")
    }) AS data,
    COLLECT(DISTINCT samp.sample_tumor_status) AS sample_tumor_status
WHERE
    (size(inputs.sample_tumor_statuses) = 0 OR size(apoc.coll.intersection(inputs.sample_tumor_statuses, sample_tumor_status)) > 0)
OPTIONAL MATCH (f)-->(:sample)-->(p:participant)
WITH inputs, na, f, 
    apoc.map.merge(data, {
        participant_id: coalesce(apoc.text.join(collect(distinct p.participant_id), ", "), "")
    }) AS data,
    COLLECT(DISTINCT p.participant_id) AS participant_id,
    COLLECT(DISTINCT p.gender) AS gender
WHERE 
    (size(inputs.subject_ids) = 0 OR size(apoc.coll.intersection(inputs.subject_ids, participant_id)) > 0) AND
    (size(inputs.genders) = 0 OR size(apoc.coll.intersection(inputs.genders, gender)) > 0)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH inputs, na, f, data,
    COLLECT(DISTINCT g.library_strategy) AS library_strategy,
    COLLECT(DISTINCT g.library_source) AS library_source,
    COLLECT(DISTINCT g.library_selection) AS library_selection,
    COLLECT(DISTINCT g.library_layout) AS library_layout,
    COLLECT(DISTINCT g.platform) AS platform,
    COLLECT(DISTINCT g.instrument_model) AS instrument_models,
    COLLECT(DISTINCT g.reference_genome_assembly) AS reference_genome_assembly
WHERE
    (size(inputs.library_strategies) = 0 OR size(apoc.coll.intersection(inputs.library_strategies, library_strategy)) > 0) AND
    (size(inputs.library_sources) = 0 OR size(apoc.coll.intersection(inputs.library_sources, library_source)) > 0) AND 
    (size(inputs.library_selections) = 0 OR size(apoc.coll.intersection(inputs.library_selections, library_selection)) > 0) AND
    (size(inputs.library_layouts) = 0 OR size(apoc.coll.intersection(inputs.library_layouts, library_layout)) > 0) AND
    (size(inputs.platforms) = 0 OR size(apoc.coll.intersection(inputs.platforms, platform)) > 0) AND
    (size(inputs.instrument_models) = 0 OR size(apoc.coll.intersection(inputs.instrument_models, instrument_models)) > 0) AND
    (size(inputs.reference_genome_assemblies) = 0 OR size(apoc.coll.intersection(inputs.reference_genome_assemblies, reference_genome_assembly)) > 0)
OPTIONAL MATCH (f)-->(:sample)-->(:participant)<--(diag:diagnosis)
WITH inputs, na, f, data,
    COLLECT(DISTINCT diag.primary_diagnosis) AS primary_diagnosis
WHERE
    (size(inputs.primary_diagnoses) = 0 OR size(apoc.coll.intersection(inputs.primary_diagnoses, primary_diagnosis)) > 0)
RETURN 
    data.file_name AS `File Name`,
    data.study_name AS `Study Name`,
    data.accession AS `Accession`,
    data.participant_id AS `Participant Id`,
    data.sample_id AS `Sample Id`,
    data.file_type AS `File Type`
ORDER BY `File Name`
LIMIT 100
'@

$ws.Range("B2").Value = $participantQuery
$ws.Range("B3").Value = $sampleQuery
$ws.Range("B4").Value = $fileQuery

$wb.Save()
